$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.10%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.65%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.116"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.42%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05665"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.86%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.472"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.8240"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.46%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8474"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.29%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.009998"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1,565.54%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1328"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.28%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.06994"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.56%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02883"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.06%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09394"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.00%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001510"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.22%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04128"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-11.89%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006224"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.22%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-1.56%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.83%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.311"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'9.08%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.03150"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.10%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1255"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.86%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.563"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-5.06%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.00%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-2.37%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004453"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-3.79%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009798"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'2.18%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'3.57%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03673"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.11%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-1.70%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1052"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.53%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-7.91%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009687"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.74%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'0.61%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.10%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-36.81%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002572"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'24.69%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
